$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1939.5
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H117").Value = 48318.668
$ws.Range("J117").Value = 48318.668
$ws.Range("L117").Value = 48318.668
$ws.Range("N117").Value = -57496.668
$ws.Range("H124").Value = 47014.2
$ws.Range("J124").Value = 47014.2
$ws.Range("L124").Value = 47014.2
$ws.Range("N124").Value = -56834.2
$ws.Range("H126").Value = 46768.8
$ws.Range("J126").Value = 46768.8
$ws.Range("L126").Value = 46768.8
$ws.Range("N126").Value = -56648.8
$ws.Range("H127").Value = 1189
$ws.Range("I127").Value = 520.8182
$ws.Range("J127").Value = 1924
$ws.Range("K127").Value = 1562.4546
$ws.Range("L127").Value = 5772
$ws.Range("M127").Value = 3397.5454
$ws.Range("N127").Value = -15692
$ws.Range("H130").Value = 48274.4
$ws.Range("J130").Value = 48274.4
$ws.Range("L130").Value = 48274.4
$ws.Range("N130").Value = -58314.4
$ws.Range("H131").Value = 2761.6365
$ws.Range("I131").Value = 2111.5
$ws.Range("K131").Value = 6334.5
$ws.Range("M131").Value = -1294.5
$ws.Range("H141").Value = 6776.25
$ws.Range("J141").Value = 10105
$ws.Range("L141").Value = 30315
$ws.Range("N141").Value = -40675

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 33000
$ws.Range("J75").Value = 33000
$ws.Range("L75").Value = 33000
$ws.Range("N75").Value = -34748
$ws.Range("H78").Value = 33000
$ws.Range("J78").Value = 33000
$ws.Range("L78").Value = 99000
$ws.Range("N78").Value = -107736
$ws.Range("H80").Value = 53323
$ws.Range("J80").Value = 53323
$ws.Range("L80").Value = 53323
$ws.Range("N80").Value = -55319
$ws.Range("H83").Value = 53323
$ws.Range("J83").Value = 53323
$ws.Range("L83").Value = 159969
$ws.Range("N83").Value = -169953
$ws.Range("H105").Value = 47954
$ws.Range("J105").Value = 47954
$ws.Range("L105").Value = 47954
$ws.Range("N105").Value = -54942
$ws.Range("H107").Value = 36272.332
$ws.Range("J107").Value = 36272.332
$ws.Range("L107").Value = 36272.332
$ws.Range("N107").Value = -43952.332
$ws.Range("H109").Value = 43977.75
$ws.Range("J109").Value = 43977.75
$ws.Range("L109").Value = 43977.75
$ws.Range("N109").Value = -46751.75
$ws.Range("H117").Value = 48367.8
$ws.Range("J117").Value = 48367.8
$ws.Range("L117").Value = 48367.8
$ws.Range("N117").Value = -57545.8
$ws.Range("H118").Value = 49582
$ws.Range("J118").Value = 49582
$ws.Range("L118").Value = 49582
$ws.Range("N118").Value = -52896
$ws.Range("H119").Value = 52640
$ws.Range("J119").Value = 52640
$ws.Range("L119").Value = 52640
$ws.Range("N119").Value = -62316

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 32191.143
$ws.Range("J104").Value = 32191.143
$ws.Range("L104").Value = 32191.143
$ws.Range("N104").Value = -37433.143
$ws.Range("H109").Value = 28715.445
$ws.Range("J109").Value = 28715.445
$ws.Range("L109").Value = 28715.445
$ws.Range("N109").Value = -30795.445
$ws.Range("H111").Value = 47279.332
$ws.Range("J111").Value = 47279.332
$ws.Range("L111").Value = 47279.332
$ws.Range("N111").Value = -55459.332
$ws.Range("H115").Value = 30790.2
$ws.Range("J115").Value = 30790.2
$ws.Range("L115").Value = 30790.2
$ws.Range("N115").Value = -33140.2
$ws.Range("H116").Value = 47832
$ws.Range("J116").Value = 47832
$ws.Range("L116").Value = 47832
$ws.Range("N116").Value = -57010
$ws.Range("H120").Value = 30107.166
$ws.Range("J120").Value = 30107.166
$ws.Range("L120").Value = 30107.166
$ws.Range("N120").Value = -37365.166

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 580.85
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 786.9286
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 2360.7858
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -2584.7858
$ws.Range("H113").Value = 2761.375
$ws.Range("I113").Value = 4363.037
$ws.Range("J113").Value = 702.0952
$ws.Range("K113").Value = 13089.111
$ws.Range("L113").Value = 2106.2856
$ws.Range("M113").Value = -10919.111
$ws.Range("N113").Value = -6446.2856

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 44960.668
$ws.Range("J104").Value = 44960.668
$ws.Range("L104").Value = 44960.668
$ws.Range("N104").Value = -51948.668
$ws.Range("H105").Value = 42791.4
$ws.Range("J105").Value = 42791.4
$ws.Range("L105").Value = 42791.4
$ws.Range("N105").Value = -49779.4
$ws.Range("H116").Value = 38967.145
$ws.Range("J116").Value = 38967.145
$ws.Range("L116").Value = 38967.145
$ws.Range("N116").Value = -48145.145
$ws.Range("H118").Value = 38256.332
$ws.Range("J118").Value = 38256.332
$ws.Range("L118").Value = 38256.332
$ws.Range("N118").Value = -41570.332
$ws.Range("H120").Value = 33531.332
$ws.Range("J120").Value = 33531.332
$ws.Range("L120").Value = 33531.332
$ws.Range("N120").Value = -43207.332

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 40002924
$ws.Range("I7").Value = 62502476
$ws.Range("J7").Value = 3722.7778
$ws.Range("K7").Value = 62502476
$ws.Range("L7").Value = 3722.7778
$ws.Range("M7").Value = -62502364
$ws.Range("N7").Value = -3946.7778
$ws.Range("H17").Value = 515
$ws.Range("I17").Value = 515
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 515
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -345
$ws.Range("N17").ClearContents()
$ws.Range("H30").Value = 2548.3333
$ws.Range("I30").Value = 2548.3333
$ws.Range("K30").Value = 2548.3333
$ws.Range("M30").Value = -2440.3333
$ws.Range("H55").Value = 518.6667
$ws.Range("I55").Value = 574.6667
$ws.Range("J55").Value = 462.66666
$ws.Range("K55").Value = 574.6667
$ws.Range("L55").Value = 462.66666
$ws.Range("M55").Value = -401.6667
$ws.Range("N55").Value = -808.66666
$ws.Range("H110").Value = 45590
$ws.Range("J110").Value = 45590
$ws.Range("L110").Value = 45590
$ws.Range("N110").Value = -53770
$ws.Range("H126").Value = 40002924
$ws.Range("I126").Value = 62502476
$ws.Range("J126").Value = 3722.7778
$ws.Range("K126").Value = 187507428
$ws.Range("L126").Value = 11168.3334
$ws.Range("M126").Value = -187504958
$ws.Range("N126").Value = -16108.3334
$ws.Range("H130").Value = 37802.332
$ws.Range("J130").Value = 37802.332
$ws.Range("L130").Value = 37802.332
$ws.Range("N130").Value = -47842.332

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64212
$ws.Range("J46").Value = 64212
$ws.Range("L46").Value = 64212
$ws.Range("N46").Value = -64674
$ws.Range("H105").Value = 50544
$ws.Range("J105").Value = 50544
$ws.Range("L105").Value = 50544
$ws.Range("N105").Value = -57532
$ws.Range("H125").Value = 39710.2
$ws.Range("J125").Value = 39710.2
$ws.Range("L125").Value = 39710.2
$ws.Range("N125").Value = -49550.2
$ws.Range("H134").Value = 64212
$ws.Range("J134").Value = 64212
$ws.Range("L134").Value = 192636
$ws.Range("N134").Value = -197706
